$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D so existing D (Tipo) becomes E, and the
# new D column gets a header of "MAE" plus numeric values.
$ws.Range("D1").EntireColumn.Insert()

# New header for column D
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# New MAE values for rows 2-4
$ws.Range("D2").Value = 0.1627423106237687
$ws.Range("D3").Value = 0.2020842214174208
$ws.Range("D4").Value = 0.2389146089630678

# Updated MSE (column B) values reflecting the re-run model
$ws.Range("B2").Value = 0.04937411793901165
$ws.Range("B3").Value = 0.07385238581709302
$ws.Range("B4").Value = 0.08729857911332642

# Updated R2 (column C) values reflecting the re-run model
$ws.Range("C2").Value = 0.9985479914260672
$ws.Range("C3").Value = 0.9992994770237257
$ws.Range("C4").Value = 0.9988180311265089
